$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the hour values for "Investigación y documentación clasificación" (E6)
# and "Clasificación k-means" (E7). The TOTAL formula in E15 (SUM(E3:E14))
# will recalculate automatically from 101.5 to 111.
$ws.Range("E6").Value = 19
$ws.Range("E7").Value = 42

# Move the active selection from E7 to E8, matching the saved cursor position.
$ws.Activate()
$ws.Range("E8").Select()
